# Duplicate the "TimeTracking" sheet, put the copy first, rename it
# "TimeTracking (2)", update its data/table to the new 2023 layout, and
# leave the original "TimeTracking" sheet (now second) intact apart from
# no longer being the active/selected tab (handled automatically by Copy).

$wb = $excel.ActiveWorkbook
$orig = $wb.Worksheets.Item("TimeTracking")

# Copy before itself -> new sheet becomes sheet 1, named "TimeTracking (2)",
# and becomes the active sheet; the original sheet loses tabSelected.
$orig.Copy($orig)

$ws = $wb.Worksheets.Item(1)

# Remove the table that got carried over conceptually (the copy does not
# actually receive a ListObject, but make sure in case it does).
if ($ws.ListObjects.Count -gt 0) {
    $ws.ListObjects.Item(1).Unlist()
}

# ---- Write the new header row ----
$headers = @("Line number","Date","Client","Project","Activity","chargeable","Tasks","Description","Start","End")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# ---- Write the new data rows (row 2 .. row 11) ----
# Columns: A Line number, B Date, C Client, D Project, E Activity,
#          F chargeable, G Tasks, H Description, I Start, J End
$data = @(
    @(1,  45100, "TimeKex", "Project 1 (TimeKex)", "TimeKex doing",    "yes", "perform tasks", $null,                  0.33333333333333331, 0.70833333333333337),
    @(2,  45101, "TimeKex", "Project 1 (TimeKex)", "TimeKex doing",    "yes", "perform tasks", $null,                  0.33333333333333331, 0.70833333333333337),
    @(3,  45102, "TimeKex", "Project 1 (TimeKex)", "TimeKex doing",    "yes", "perform tasks", $null,                  0.33333333333333331, 0.54166666666666663),
    @(4,  45102, "TimeKex", "Project 1 (TimeKex)", "TimeKex doing",    "yes", "perform tasks", $null,                  0.625,                0.70833333333333337),
    @(5,  45103, "TimeKex", "Project 1 (TimeKex)", "TimeKex learning", "yes", "learn stuff",    "will not be billable", 0.33333333333333331, 0.41666666666666669),
    @(6,  45103, "TimeKex", "Project 1 (TimeKex)", "TimeKex doing",    "yes", "perform tasks", "will be billable",     0.4375,               0.5),
    @(7,  45103, "TimeKex", "Project 1 (TimeKex)", "TimeKex doing",    "yes", "perform tasks", "will be billable",     0.54513888888888895, 0.72222222222222221),
    @(8,  45104, "TimeKex", "Project 1 (TimeKex)", "TimeKex doing",    "yes", "perform tasks", "will be billable",     0.33333333333333331, 0.4201388888888889),
    @(9,  45104, "TimeKex", "Project 1 (TimeKex)", "TimeKex doing",    "yes", "perform tasks", "will be billable",     0.4201388888888889,  0.5),
    @(10, 45104, "TimeKex", "Project 1 (TimeKex)", "TimeKex doing",    "no",  "perform tasks", "will not be billable", 0.53819444444444442, 0.65625)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($r + 2, $c + 1).Value = $val
        }
    }
}

# Clear any leftover cells outside the new A1:J11 extent (old sheet had
# K column and 14 rows).
$ws.Range("K1:K14").Clear()
$ws.Range("A12:J14").Clear()

# ---- Number formats ----
# Date column -> yyyy-mm-dd ; Start/End columns keep h:mm (inherited).
$ws.Range("B2:B11").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("I2:J11").NumberFormat = "h:mm;@"

# ---- Column widths (match target layout after column removal) ----
$ws.Columns.Item(2).ColumnWidth = 14.42578125
$ws.Columns.Item(7).ColumnWidth = 14.85546875
$ws.Columns.Item(8).ColumnWidth = 20.85546875

# ---- Rebuild the table over the new extent ----
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:J11"), $null, 1)
$tbl.Name = "Tabelle22"
$tbl.TableStyle = "TableStyleMedium2"

# ---- View / selection ----
$ws.Activate()
$ws.Range("H17").Select()

Write-Output "done"
